$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "311.98"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-4.55%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-7.04%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.103"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.25%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07878"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-5.64%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.331"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.91%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.689"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-12.83%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9258"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.47%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1082"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.84%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1785"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.64%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09104"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-6.43%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04403"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.19%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.192"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-16.40%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.14%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001267"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.59%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006026"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.17%"
$ws.Range("B17").Value = "HotbitToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.004152"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-5.88%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.378"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.72%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.559"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.71%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3317"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.83%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1376"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.36%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2800"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.65%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04157"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.01%"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001242"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.70%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.79%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0002994"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.36%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02453"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-8.63%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05327"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.008020"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.13%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1357"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.66%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007580"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.29%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001997"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.45%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008205"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "4.75%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3104"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-11.40%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006761"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.06%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000754"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.35%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003432"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.65%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004120"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "16.56%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002110"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.35%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002010"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.35%"
